$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Tyrell"
$ws.Range("B2").Value = "Jones"
$ws.Range("C2").Value = "rickie.lebsack@gmail.com"
$ws.Range("D2").Value = "1b581zkh"
$ws.Range("E2").Value = "46674 Trantow Grove"
